$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "Website" (C) and "Branch" (D) columns so the old "Marks" column
# (E) shifts left into C, leaving just Name / Email / Marks.
$ws.Range("C:D").Delete()

# Header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Marks"

# Merged data rows (all three rows now share the same values)
$ws.Range("A2:A4").Value = "gggfd garg"
$ws.Range("B2:B4").Value = "gdgg@fc.in"
$ws.Range("C2:C4").Value = 343
